$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Wipe the old layout (rows 1-7, cols A-E) and rebuild it fresh so that no
# shift/copy re-serializes the existing floating point numbers (the old
# "Adult Survival" row also disappears this way).
$ws.Range("A1:E7").EntireRow.Delete()

# Header row: new LH_Stage / Month columns in front of the existing stats.
$ws.Range("A1").Value = "LH_Stage "
$ws.Range("B1").Value = "Month "
$ws.Range("C1").Value = "MinTemp"
$ws.Range("D1").Value = "MaxTemp"
$ws.Range("E1").Value = "AvgTemp"
$ws.Range("F1").Value = "SD"

# Eggs
$ws.Range("A2").Value = "Eggs"
$ws.Range("B2").Value = "01, 02, 03, 11, 12"
$ws.Range("C2").Value = -3.197
$ws.Range("D2").Value = 11.528
$ws.Range("E2").Value = 4.0784921185154497
$ws.Range("F2").Value = 2.4001288794104401

# Alevin
$ws.Range("A3").Value = "Alevin"
$ws.Range("B3").Value = "02, 03, 04"
$ws.Range("C3").Value = -1.456
$ws.Range("D3").Value = 17.95
$ws.Range("E3").Value = 4.6695105337078697
$ws.Range("F3").Value = 3.02052130260072

# YOY
$ws.Range("A4").Value = "YOY"
$ws.Range("B4").Value = "06, 07, 09, 09, 10"
$ws.Range("C4").Value = 5.655
$ws.Range("D4").Value = 20.138000000000002
$ws.Range("E4").Value = 13.7522661626928
$ws.Range("F4").Value = 2.7718664686113899

# Spawning
$ws.Range("A5").Value = "Spawning"
$ws.Range("B5").Value = "09, 10, 11"
$ws.Range("C5").Value = 2.837
$ws.Range("D5").Value = 18.901
$ws.Range("E5").Value = 10.149747052426701
$ws.Range("F5").Value = 3.1867817403018002

# New "Month" column is a bit wider than the default.
$ws.Columns.Item(2).ColumnWidth = 14.5

# Selection left as it was when the author last saved.
$ws.Range("B10:B14").Select() | Out-Null
